# Term ValueSet KLGeneralInformationFSIII - bump to version 2.0.0, refresh date/contact,
# and add a second "Include from FSIII" sheet (copy of the first).

$wb = $excel.ActiveWorkbook

# --- 1. Update Metadata sheet values ---------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- 2. Duplicate "Include from FSIII" sheet, placed right after it --------
$src = $wb.Worksheets.Item("Include from FSIII")
$src.Copy([System.Reflection.Missing]::Value, $src)

$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "Include from FSIII 2"
